$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 758, shifting the existing weekly
# Primera/Segunda pairs (rows 758-809) down to rows 760-811.
$ws.Range("A758:A759").EntireRow.Insert()

# Populate the two newly inserted rows with the new weekly observation
# (date 44826) that now heads this block.

# Row 758 - "Primera" quality
$ws.Cells.Item(758,1).Value2  = 6
$ws.Cells.Item(758,2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(758,3).Value2  = "Metropolitana"
$ws.Cells.Item(758,4).Value2  = 44826
$ws.Cells.Item(758,5).Value2  = 13
$ws.Cells.Item(758,6).Value2  = 100112017
$ws.Cells.Item(758,7).Value2  = "Apio"
$ws.Cells.Item(758,8).Value2  = "Americana (o)"
$ws.Cells.Item(758,9).Value2  = "Primera"
$ws.Cells.Item(758,10).Value2 = 1530
$ws.Cells.Item(758,11).Value2 = 7000
$ws.Cells.Item(758,12).Value2 = 8000
$ws.Cells.Item(758,13).Value2 = 7444
$ws.Cells.Item(758,14).Value2 = "`$/docena de matas"
$ws.Cells.Item(758,15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(758,16).Value2 = 1241
$ws.Cells.Item(758,17).Value2 = 6
$ws.Cells.Item(758,18).Value2 = "Hortaliza"

# Row 759 - "Segunda" quality
$ws.Cells.Item(759,1).Value2  = 6
$ws.Cells.Item(759,2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(759,3).Value2  = "Metropolitana"
$ws.Cells.Item(759,4).Value2  = 44826
$ws.Cells.Item(759,5).Value2  = 13
$ws.Cells.Item(759,6).Value2  = 100112017
$ws.Cells.Item(759,7).Value2  = "Apio"
$ws.Cells.Item(759,8).Value2  = "Americana (o)"
$ws.Cells.Item(759,9).Value2  = "Segunda"
$ws.Cells.Item(759,10).Value2 = 470
$ws.Cells.Item(759,11).Value2 = 6000
$ws.Cells.Item(759,12).Value2 = 6000
$ws.Cells.Item(759,13).Value2 = 6000
$ws.Cells.Item(759,14).Value2 = "`$/docena de matas"
$ws.Cells.Item(759,15).Value2 = "Región de Coquimbo"
$ws.Cells.Item(759,16).Value2 = 1000
$ws.Cells.Item(759,17).Value2 = 6
$ws.Cells.Item(759,18).Value2 = "Hortaliza"
